# Timeline.xlsx update:
#  - Row 6: remove the "Currently in Progress" note in G6 and instead
#    record an Actual-finish date of 6-Oct-2017 in F6.
#  - Row 7: record an Actual-start date of 6-Oct-2017 in E7.
#  - Update the view's selection to F7 (the last cell touched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("Create Solution project scratch"): set Actual/Finished (F6) and
# clear the old "Currently in Progress" note from G6.
$ws.Range("F6").Value = 43014
$ws.Range("F6").NumberFormat = "d-mmm-yy"
$ws.Range("G6").Clear()

# Row 7 ("Prepare development libraries"): set Actual/Start (E7).
$ws.Range("E7").Value = 43014
$ws.Range("E7").NumberFormat = "d-mmm-yy"

# Reflect the new selection / scroll position left in the sheet view.
$ws.Range("F7").Select()
